# Auto-generated Excel COM-interop script
# Applies scheduled-runner market-data refresh to the Kujata Profits workbook.
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) per leve row
# across all job sheets, matching upstream API refresh output. A handful of cells
# whose upstream value is no longer present (e.g. no HQ price data) are cleared so
# the cell itself is removed rather than left holding a stale number.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 185.9375
$ws.Range("I33").Value = 128.35715
$ws.Range("K33").Value = 128.35715
$ws.Range("M33").Value = 100.64285
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("M60").ClearContents()
$ws.Range("H64").Value = 4345.385
$ws.Range("I64").Value = 4299
$ws.Range("K64").Value = 4299
$ws.Range("M64").Value = -4051
$ws.Range("H67").Value = 4345.385
$ws.Range("I67").Value = 4299
$ws.Range("K67").Value = 4299
$ws.Range("M67").Value = -3441
$ws.Range("H86").Value = 5568.5713
$ws.Range("I86").Value = 5996.6665
$ws.Range("J86").Value = 3000
$ws.Range("K86").Value = 5996.6665
$ws.Range("L86").Value = 3000
$ws.Range("M86").Value = -4873.6665
$ws.Range("N86").Value = -5246
$ws.Range("H89").Value = 5568.5713
$ws.Range("I89").Value = 5996.6665
$ws.Range("J89").Value = 3000
$ws.Range("K89").Value = 29983.3325
$ws.Range("L89").Value = 15000
$ws.Range("M89").Value = -24367.3325
$ws.Range("N89").Value = -26232
$ws.Range("H132").Value = 6178115.5
$ws.Range("I132").Value = 9528976
$ws.Range("J132").Value = 5478.421
$ws.Range("K132").Value = 28586928
$ws.Range("L132").Value = 16435.263
$ws.Range("M132").Value = -28584398
$ws.Range("N132").Value = -21495.263
$ws.Range("H137").Value = 1110.6177
$ws.Range("I137").Value = 813.75
$ws.Range("J137").Value = 1534.7142
$ws.Range("K137").Value = 2441.25
$ws.Range("L137").Value = 4604.142599999999
$ws.Range("M137").Value = 108.75
$ws.Range("N137").Value = -9704.142599999999
$ws.Range("H138").Value = 1249.56
$ws.Range("I138").Value = 610.6177
$ws.Range("J138").Value = 1578.7122
$ws.Range("K138").Value = 1831.8531
$ws.Range("L138").Value = 4736.1366
$ws.Range("M138").Value = 3308.1469
$ws.Range("N138").Value = -15016.1366
$ws.Range("H141").Value = 717.375
$ws.Range("I141").Value = 631.86664
$ws.Range("K141").Value = 1895.59992
$ws.Range("M141").Value = 3284.40008

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4147.636
$ws.Range("I32").Value = 3626.5789
$ws.Range("K32").Value = 3626.5789
$ws.Range("M32").Value = -3339.5789
$ws.Range("H61").Value = 29412982
$ws.Range("I61").Value = 37037868
$ws.Range("J61").Value = 2714.2856
$ws.Range("K61").Value = 37037868
$ws.Range("L61").Value = 2714.2856
$ws.Range("M61").Value = -37037656
$ws.Range("N61").Value = -3138.2856
$ws.Range("H74").Value = 1111.8
$ws.Range("I74").Value = 845.6667
$ws.Range("J74").Value = 3507
$ws.Range("K74").Value = 845.6667
$ws.Range("L74").Value = 3507
$ws.Range("M74").Value = 28.33330000000001
$ws.Range("N74").Value = -5255
$ws.Range("H77").Value = 1111.8
$ws.Range("I77").Value = 845.6667
$ws.Range("J77").Value = 3507
$ws.Range("K77").Value = 4228.3335
$ws.Range("L77").Value = 17535
$ws.Range("M77").Value = 139.6665000000003
$ws.Range("N77").Value = -26271
$ws.Range("H136").Value = 29412982
$ws.Range("I136").Value = 37037868
$ws.Range("J136").Value = 2714.2856
$ws.Range("K136").Value = 111113604
$ws.Range("L136").Value = 8142.8568
$ws.Range("M136").Value = -111111054
$ws.Range("N136").Value = -13242.8568

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H43").Value = 192342
$ws.Range("J43").Value = 192342
$ws.Range("L43").Value = 192342
$ws.Range("N43").Value = -192704
$ws.Range("H107").Value = 1070.5312
$ws.Range("I107").Value = 849.85187
$ws.Range("K107").Value = 849.85187
$ws.Range("M107").Value = 1070.14813
$ws.Range("H134").Value = 5670.04
$ws.Range("I134").Value = 1397.1364
$ws.Range("J134").Value = 37004.668
$ws.Range("K134").Value = 4191.4092
$ws.Range("L134").Value = 111014.004
$ws.Range("M134").Value = -1656.4092
$ws.Range("N134").Value = -116084.004

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1853.3334
$ws.Range("I31").Value = 1853.3334
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 1853.3334
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -1558.3334
$ws.Range("N31").ClearContents()
$ws.Range("H34").Value = 1853.3334
$ws.Range("I34").Value = 1853.3334
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 1853.3334
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -1651.3334
$ws.Range("N34").ClearContents()
$ws.Range("H58").Value = 857.9184
$ws.Range("I58").Value = 783.9535
$ws.Range("K58").Value = 783.9535
$ws.Range("M58").Value = -580.9535
$ws.Range("H62").Value = 5156097
$ws.Range("I62").Value = 2377.361
$ws.Range("J62").Value = 67000732
$ws.Range("K62").Value = 2377.361
$ws.Range("L62").Value = 67000732
$ws.Range("M62").Value = -1753.361
$ws.Range("N62").Value = -67001980
$ws.Range("H65").Value = 5156097
$ws.Range("I65").Value = 2377.361
$ws.Range("J65").Value = 67000732
$ws.Range("K65").Value = 11886.805
$ws.Range("L65").Value = 335003660
$ws.Range("M65").Value = -8766.805
$ws.Range("N65").Value = -335009900
$ws.Range("H134").Value = 15626187
$ws.Range("I134").Value = 1122.091
$ws.Range("K134").Value = 3366.273
$ws.Range("M134").Value = -831.2729999999997
$ws.Range("H136").Value = 857.9184
$ws.Range("I136").Value = 783.9535
$ws.Range("K136").Value = 2351.8605
$ws.Range("M136").Value = 198.1395000000002

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 662.2
$ws.Range("J113").Value = 689.4666999999999
$ws.Range("L113").Value = 2068.4001
$ws.Range("N113").Value = -6408.4001
$ws.Range("H131").Value = 32262344
$ws.Range("I131").Value = 111111660
$ws.Range("J131").Value = 5802.091
$ws.Range("K131").Value = 333334980
$ws.Range("L131").Value = 17406.273
$ws.Range("M131").Value = -333329940
$ws.Range("N131").Value = -27486.273

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()
$ws.Range("H132").Value = 2623.2144
$ws.Range("I132").Value = 2421.1875
$ws.Range("J132").Value = 2892.5833
$ws.Range("K132").Value = 7263.5625
$ws.Range("L132").Value = 8677.749899999999
$ws.Range("M132").Value = -4733.5625
$ws.Range("N132").Value = -13737.7499

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 18338.322
$ws.Range("I132").Value = 809.4737
$ws.Range("J132").Value = 50057.19
$ws.Range("K132").Value = 2428.4211
$ws.Range("L132").Value = 150171.57
$ws.Range("M132").Value = 101.5789
$ws.Range("N132").Value = -155231.57
$ws.Range("H136").Value = 1817.9
$ws.Range("I136").Value = 1650.4706
$ws.Range("J136").Value = 2766.6667
$ws.Range("K136").Value = 4951.4118
$ws.Range("L136").Value = 8300.000100000001
$ws.Range("M136").Value = -2401.4118
$ws.Range("N136").Value = -13400.0001

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 675.9655
$ws.Range("I136").Value = 617.3043
$ws.Range("J136").Value = 900.8333
$ws.Range("K136").Value = 1851.9129
$ws.Range("L136").Value = 2702.4999
$ws.Range("M136").Value = 698.0871
$ws.Range("N136").Value = -7802.4999
